$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns:
#  - one before old column F (Pri. F) -> becomes new column F ("Sample UFID, F")
#  - one before old column M (Pri. M, which after the first insert sits at N) -> becomes new column N ("Sample UFID, M")
$ws.Columns("F").EntireColumn.Insert()
$ws.Columns("N").EntireColumn.Insert()

# Column-insert shifts cell values/styles automatically, but header-row cell
# comments stay anchored to their old column letters, so migrate them by hand
# -- process from the right-most source column down to the left-most so a
# destination is never clobbered before it has been read.
$moves = @(
  @("V","X"), @("U","W"), @("T","V"), @("S","U"), @("R","T"), @("Q","S"),
  @("P","R"), @("O","Q"), @("N","P"), @("M","O"),
  @("K","L"), @("J","K"), @("I","J"), @("H","I"), @("G","H"), @("F","G")
)
foreach ($move in $moves) {
  $srcRef = $move[0] + "3"
  $dstRef = $move[1] + "3"
  $src = $ws.Range($srcRef)
  $comment = $src.Comment
  if ($comment -ne $null) {
    $text = $comment.Text()
    $author = $comment.Author
    $comment.Delete()
    $newComment = $ws.Range($dstRef).AddComment($text)
    $newComment.Author = $author
  }
}

# New header labels
$ws.Range("F3").Value = "Sample UFID, F"
$ws.Range("N3").Value = "Sample UFID, M"

# New header-cell comments (matching the ones already used for the PIT columns)
$noteText = "Optional, must match exisiting sample number.`nPIT tag field must be blank."
$ws.Range("F3").AddComment($noteText)
$ws.Range("N3").AddComment($noteText)

# Resize the two new columns to fit their header text, like the rest of the header columns
$ws.Columns("F").EntireColumn.AutoFit()
$ws.Columns("N").EntireColumn.AutoFit()

# Restore the cursor/selection position left behind by the edit
$ws.Range("G18").Select()
